# "aanpassing documentatie op basis van besluit objecttoestand"
# Sort the hoofdgroepen table alphabetically on column C (hoofdgroep),
# turn the range into an AutoFilter table, remember the sort as the
# filter's sort state, resize columns A and C to fit their content,
# and update the active view (scroll position / selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("A1:D33")
$sortKey   = $ws.Range("C1:C33")

# Sort rows 2-33 ascending by column C, keeping the header row (row 1) in place.
$dataRange.Sort(
    $sortKey,
    [Microsoft.Office.Interop.Excel.XlSortOrder]::xlAscending,
    $null,
    $null,
    [Microsoft.Office.Interop.Excel.XlSortOrder]::xlAscending,
    $null,
    [Microsoft.Office.Interop.Excel.XlSortOrder]::xlAscending,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)

# Turn the sorted range into a filtered table (adds autoFilter + the
# hidden _xlnm._FilterDatabase defined name).
$dataRange.AutoFilter()

$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "='5.01_hoofdgroepen'!`$A`$1:`$D`$33")
$filterName.Visible = $false

# Fit column A and column C widths to their (now longest) content.
$ws.Columns.Item(1).ColumnWidth = 70.3
$ws.Columns.Item(3).ColumnWidth = 28.3

# Update the window view: scroll so row 16 is at the top, and make C10 the
# active / selected cell.
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
$ws.Range("C10").Select()
